$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 36 (shifts the existing rows 36-60 down to 37-61)
# and populate it with a new expression/answer pair.
$ws.Rows.Item(36).Insert()
$ws.Range("A36").Value = "(25 - 3) * 2"
$ws.Range("B36").Value = 44

# Rows 37-55 (formerly 36-54) retained their original expression/answer pairs
# automatically via the row insert/shift above - no further action needed there.

# Row 56 (formerly row 55, "(8 - 3) * 4" / 20) becomes a new expression with no answer yet.
$ws.Range("A56").Value = "20 ÷ 50 + 8"
$ws.Range("B56").ClearContents()

# Row 57 (formerly row 56, "200 ÷ 50 + 8" / 12) becomes a new expression with no answer yet.
$ws.Range("A57").Value = "(18 - 3) * 4"
$ws.Range("B57").ClearContents()

# Row 58 (formerly row 57, "(10 * 2 ) - (3 * 4)" / 8) is unchanged.

# Row 59 (formerly row 58, "(9 + 1) - (9 - 1)" / 80) gets reworded; the answer stays the same.
$ws.Range("A59").Value = "(9 + 1) * (9 - 1)"

# Rows 60-61 (formerly rows 59-60, "4 * (3 + 5) - 7" / 25 and "8 * (7 - 3)" / 40) are unchanged.

# Leave the selection where the author ended up working.
$ws.Range("C37").Select()
